# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with new Binance rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.9 = 6889.88 pesos`n✅ 6889.88 pesos = 1.89 = 936.95 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas sheet: update N10/O10/N12/O12 rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 525.5
$ws2.Range("O10").Value = 3620.63
$ws2.Range("N12").Value = 3640
$ws2.Range("O12").Value = 495.001
